# Actualizacion desde MV -datos-
# Appends the five newest "Diaria" observations (27-09-2021 .. 01-10-2021)
# to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Fecha = "27-09-2021"; Uno = 3.22; Tres = 3.19; Cinco = 3.28 },
    @{ Fecha = "28-09-2021"; Uno = 3.23; Tres = 3.2;  Cinco = 3.27 },
    @{ Fecha = "29-09-2021"; Uno = 3.22; Tres = 3.24; Cinco = 3.33 },
    @{ Fecha = "30-09-2021"; Uno = 3.21; Tres = 3.17; Cinco = 3.24 },
    @{ Fecha = "01-10-2021"; Uno = 3.4;  Tres = 3.18; Cinco = 3.23 }
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$scratch = $ws.Range("ZZ1")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # The "Serie" column holds dd-mm-yyyy text. Most of these values are
    # unambiguous as text (day > 12) and can be typed in directly, but a
    # few (e.g. 01-10-2021) look like a valid date to Excel's type
    # inference and would otherwise be auto-converted to a date serial.
    # Build the label as a text formula result first, then paste its
    # *value* into the real cell so it lands as plain text, exactly like
    # its neighbours - Excel's paste-values path doesn't re-run the
    # "does this look like a date" inference that live typing does.
    $scratch.Formula = '="' + $row.Fecha + '"'
    $scratch.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $row.Uno
    $ws.Cells.Item($r, 3).Value = $row.Tres
    $ws.Cells.Item($r, 4).Value = $row.Cinco
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
